# Automatische test-sync: 2025-06-19 12:30:10
#
# Appends a new "Afmelding nieuwsbrief" (newsletter opt-out) log entry to the
# Logs sheet (row 10) and bumps the "Afmelding" tally on the Dashboard sheet
# from 1 to 2.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(10, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item(10, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(10, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item(10, 4).Value = "Afmelding"
# Column E (Antwoord) stays empty for this row, same as the source record.
$logs.Cells.Item(10, 6).Value = "2025-06-19 12:28:11"
$logs.Cells.Item(10, 7).Value = "Nee"

# Extend the existing conditional-formatting rules (Categorie / Beantwoord)
# so they keep covering the sheet's used range, now through row 10.
$catRules = $logs.Range("D2:D9").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D10"))
}

$answeredRules = $logs.Range("G2:G9").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G10"))
}

# --- Dashboard sheet: update the "Afmelding" count -------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 2
